$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# 1) Structural edits: remove the now-unused blank rows.
#    - Rows 10:11 (leftover blank "Bruna Filgueiras" continuation rows)
#      are removed, shifting everything below up by 2.
#    - The blank spacer row that used to be row 31 (now row 29 after the
#      first shift) is also removed, shifting everything below up by 1
#      more (net -3 from the original row 32 onward).
# -----------------------------------------------------------------
$ws.Rows("10:11").Delete()
$ws.Rows("29").Delete()

# -----------------------------------------------------------------
# 2) Update the title cell (now referencing a different shared string)
# -----------------------------------------------------------------
$ws.Range("D1").Value = "Notas: AV1 - Projetos Python IoT - 01/05/2023"

# -----------------------------------------------------------------
# 3) Cell content / value updates on the already-shifted grid
# -----------------------------------------------------------------

# Team 1 - Franklin Souza
$ws.Range("E5").Value = 7
$ws.Range("F5").Value = 0
$ws.Range("H5").Value = "Web; Github; Lista de Tarefas; Não edita; não consulta; https://projetowebpython.onrender.com/"

# Team 2 - Bruna Filgueiras
$ws.Range("E9").Value = 7
$ws.Range("F9").Value = 1
$ws.Range("H9").Value = "Web; Github (não); Leitor de Scripts; não remove; não guarda no BD; https://filgueiras7-leitorscript-app-x4nmz3.streamlit.app/"

# Team 3 - Angelo do Nascimento
$ws.Range("E11").Value = 9
$ws.Range("F11").Value = 0
$ws.Range("H11").Value = "Github; Web; Cadastro de Fornecedores; informar o atributo de busca; validar atributos, ex. e-mail; https://projetofullstackpython-production.up.railway.app/"

# Team 4 - Joaquim Juliao / Everton Mutti / Davi Chagas / Matheus Matos
$ws.Range("H15").Value = "Github; API para integrar produtos Shp; Time Solidário (ensino API/GIT). https://github.com/EvertonMutti/API_Shopping"
$ws.Range("E16").Value = 7.8
$ws.Range("F16").Value = 1
$ws.Range("F18").Value = 1

# Team 5 - Levi / Mercia / Flavio Castro / Vinicius Scandura
$ws.Range("E20").Value = 7
$ws.Range("F20").Value = 0
$ws.Range("H20").Value = "Github não; Universo Colaborativo; VSCode; Posts (consulta; delete; alteração); js dinâmico; sem API"

# Team 6 - Adrielle Santana / Caike Rocha / João Pamponet / Klaus Erick Maciel
$ws.Range("E25").Value = 7
$ws.Range("F25").Value = 1
$ws.Range("H25").Value = "Github; Projeto Arduino; https://github.com/JV1T0R/Arduino-project; Projeto não rodando"
$ws.Range("F27").Value = 1
$ws.Range("H29").Value = "Parabéns Klaus, envio de e-mail; Github; enviar email; Não Equipe de João Pamponet Projeto Arduino"

# Team 7 - Raudiney Andrade / David / Rafael / Vinicius
$ws.Range("E31").Value = 7
$ws.Range("F31").Value = 1
$ws.Range("H31").Value = "WEB(não); Github(Não); Vídeo MKV; BI Cotações do Ibovespa; faltou a analise dos dados; quando investir; quanto investir; consultoria financeira."

# Team 8 - Ariel Araujo
$ws.Range("E36").Value = 6
$ws.Range("F36").Value = 0

# -----------------------------------------------------------------
# 4) New rows at the bottom for the extra quiz entries (Thauan, Luiz Junior)
# -----------------------------------------------------------------
$ws.Range("C40").Value = "Quiz"
$ws.Range("D40").Value = "Thauan"
$ws.Range("F40").Value = 1

$ws.Range("C41").Value = "Quiz"
$ws.Range("D41").Value = "Luiz Junior"
$ws.Range("F41").Value = 1

# -----------------------------------------------------------------
# 5) Update the view (scroll position / active selection)
# -----------------------------------------------------------------
$ws.Range("E20").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 9
$excel.ActiveWindow.ScrollColumn = 1
